$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Remove the "Meta description: ..." paragraph that used to follow
#    the document's H1 title.
# ------------------------------------------------------------------
$metaOld = "Meta description: Explore the Amazon Rainforest in the Amazonia online slot game by Merkur. Enjoy mini-games, free spins, and unique symbols. Play for free now."
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.StartsWith("Meta description")) {
        $p.Range.Delete()
        break
    }
}

# ------------------------------------------------------------------
# 2. Insert a new bold paragraph ("Play Amazonia Slot for Free: Review
#    & Features") right before the last paragraph (the former "Prompt:"
#    paragraph).
# ------------------------------------------------------------------
$count = $d.Paragraphs.Count
$lastPara = $d.Paragraphs.Item($count)
$insertPos = $d.Range($lastPara.Range.Start, $lastPara.Range.Start)

$xmlFrag = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Play Amazonia Slot for Free: Review &amp; Features</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"></w:p>'
$insertPos.InsertXML($xmlFrag)

# InsertXML needs a second (dummy) paragraph mark to force the real
# break; remove that now-empty spacer paragraph.
$dummyIndex = $count + 1
$dummyPara = $d.Paragraphs.Item($dummyIndex)
$dummyRange = $d.Range($dummyPara.Range.Start, $dummyPara.Range.End)
$dummyRange.Delete()

# ------------------------------------------------------------------
# 3. Replace the old "Prompt: ..." image-prompt text with the meta
#    description text (keeping the paragraph's italic formatting).
# ------------------------------------------------------------------
$old = "Prompt: Create a cartoon-style feature image for Amazonia slot game. The image should depict a happy Maya warrior with glasses. The background should showcase the lush green of the Amazon rainforest. The Maya warrior should be holding a tablet or smartphone with the Amazonia game logo on it. The image should convey a fun and exciting gaming experience with the Amazonia game in a playful manner. The colours should be bright and vibrant, providing an eye-catching contrast to the green background. Please ensure that the image is in high-resolution to be used not only in the game but for promotional purposes too."
$new = "Explore the Amazon Rainforest in the Amazonia online slot game by Merkur. Enjoy mini-games, free spins, and unique symbols. Play for free now."
$d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)

Write-Host "Edit complete. Paragraph count:" $d.Paragraphs.Count
